$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The comparison_title column (G) had duplicate "IMPDH1 knockdown" titles for
# both the 1v2 (G2) and 1v3 (G3) comparisons. Row 3 actually compares
# shControl vs shIMPDH2 (see F3 = "1v3"), so its title should read
# "IMPDH2 knockdown" instead of the duplicated "IMPDH1 knockdown".
$ws.Range("G3").Value = "IMPDH2 knockdown"
